# The document has a single paragraph whose entire text is a real
# w:hyperlink (external relationship) reading "ex05: <url>". The target
# state instead spells the same visible text with a leading "e" typed
# back in front of a separately-split "x05: " run, followed by the
# hyperlink re-expressed as a classic Word field:
#   { HYPERLINK "https://www.youtube.com/watch?v=Xe_rS_yTqJQ" }
# i.e. fldChar(begin) / instrText*3 / fldChar(separate) / result-run
# (still styled with the Hyperlink character style) / fldChar(end).
# This is exactly what Word leaves behind when a hyperlink's display
# text is edited in place and the link gets converted to a field code,
# so we reconstruct it by dropping the old w:hyperlink wrapper and
# splicing in the equivalent run/field markup over the same text.

$d = $word.ActiveDocument

$url = "https://www.youtube.com/watch?v=Xe_rS_yTqJQ"
$displayText = "ex05: " + $url

# Step 1: strip the w:hyperlink wrapper. Hyperlink.Delete() unlinks it
# (drops the field/relationship) while leaving the run text + the
# direct "Hyperlink" character-style formatting (rStyle a3) in place,
# which is exactly the formatting the new result-run still needs.
if ($d.Hyperlinks.Count -gt 0) {
    $d.Hyperlinks.Item(1).Delete()
}

# Step 2: find the (now plain) run text so we can replace just that
# span, leaving the paragraph mark / pPr untouched.
$rng = $d.Content.Duplicate
$ok = $rng.Find.Execute($displayText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Step 3: rebuild that span as: "e" / "x05: " / begin / instrText*3 /
# separate / url-result (styled) / end - all sized 30/30 like the
# original runs, only the field-result run keeping rStyle "a3".
$sz = '<w:rPr><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr>'
$szLinked = '<w:rPr><w:rStyle w:val="a3"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr>'

$runs = ""
$runs += "<w:r>$sz<w:t>e</w:t></w:r>"
$runs += "<w:r>$sz<w:t xml:space=`"preserve`">x05: </w:t></w:r>"
$runs += "<w:r>$sz<w:fldChar w:fldCharType=`"begin`"/></w:r>"
$runs += "<w:r>$sz<w:instrText xml:space=`"preserve`"> HYPERLINK `"</w:instrText></w:r>"
$runs += "<w:r>$sz<w:instrText>$url</w:instrText></w:r>"
$runs += "<w:r>$sz<w:instrText xml:space=`"preserve`">`" </w:instrText></w:r>"
$runs += "<w:r>$sz<w:fldChar w:fldCharType=`"separate`"/></w:r>"
$runs += "<w:r>$szLinked<w:t>$url</w:t></w:r>"
$runs += "<w:r>$sz<w:fldChar w:fldCharType=`"end`"/></w:r>"

$xmlFrag = "<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?>" +
  "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" +
  "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" +
  "<pkg:xmlData>" +
  "<w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">" +
  "<w:body><w:p>$runs</w:p></w:body>" +
  "</w:document>" +
  "</pkg:xmlData></pkg:part></pkg:package>"

if ($ok) {
    $rng.InsertXML($xmlFrag)
} else {
    # Fallback: the whole paragraph (minus its mark) is the hyperlink text.
    $d.Paragraphs.Item(1).Range.InsertXML($xmlFrag)
}
